$d = $word.ActiveDocument

# --- 1. Paragraph 3: merge "with almost " + "all of" + " the growth..." runs
#        into a single run (drops the gramStart/gramEnd proofErr markers),
#        by re-finding & replacing the span with itself.
$find1 = $d.Paragraphs.Item(3).Range.Find
$find1.Execute("with almost all of the growth", $true, $false, $false, $false, $false, $true, 1, $false, "with almost all of the growth", 2)

# --- 2. Paragraph 4 (the trailing empty paragraph) gets two new sentences
#        about Cambodia's forest cover, added as two separate runs.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Text = "At the start of the century 41.9% of Cambodia’s land area was forested"

# Add a new paragraph right after it, and fill that with the 2nd sentence.
$d.Paragraphs.Item(4).Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$p5.Range.Text = ", and by 2012 the total forested area had been reduced by 19.8%, equating to over 1.3 million hectares (Davis et al 2015). Only 25 other countries lost more forest than Cambodia between 2000 – 2012 (Hansen et al 2013). "

# Merge the two paragraphs back together by deleting the paragraph mark that
# separates them -- this keeps the two sentences as two distinct runs
# instead of collapsing them into one run.
$p4again = $d.Paragraphs.Item(4)
$endOfP4 = $p4again.Range.End
$markRange = $d.Range($endOfP4 - 1, $endOfP4)
$markRange.Delete()
